# Analyse des performances apres la 6eme correction
# Fill in the "ETAPE 6" column (I) on each sheet by duplicating the
# "ETAPE 5" column (H) values/formatting, then update the selections to
# reflect the newly-populated column.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: "LightHouse - Portable" -------------------------------------
$ws1 = $wb.Worksheets.Item("LightHouse - Portable")
$ws1.Range("H4:H7").Copy($ws1.Range("I4:I7")) | Out-Null
$ws1.Range("H4:I7").Select() | Out-Null

# --- Sheet 2: "LightHouse - Bureau" ----------------------------------------
$ws2 = $wb.Worksheets.Item("LightHouse - Bureau")
$ws2.Range("H5:H7").Copy($ws2.Range("I5:I7")) | Out-Null
# Row 4 on this sheet got a fresh measurement (88) rather than a repeat of
# column H (92) - copy the column-H formatting pattern for row 4 (style
# already present on D4 of the same row) and then overwrite the value.
$ws2.Range("D4").Copy($ws2.Range("I4")) | Out-Null
$ws2.Range("I4").Value = 88
$ws2.Range("I6:I7").Select() | Out-Null

# --- Sheet 3: "GTmetrix - Bureau" ------------------------------------------
$ws3 = $wb.Worksheets.Item("GTmetrix - Bureau")
$ws3.Range("H4:H5").Copy($ws3.Range("I4:I5")) | Out-Null
$ws3.Range("H4:I5").Select() | Out-Null
